$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5, shifting existing rows 5:39 down to 6:40
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data record
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44558
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103006
$ws.Range("J5").Value = "Nectarín"
$ws.Range("K5").Value = "Artic Star"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19000
$ws.Range("Q5").Value = "$/bandeja 18 kilos granel"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1056
$ws.Range("T5").Value = 18
